$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 ("nav1 h1" breakpoint font sizes): tweak the 768px and 992px values,
# and drop the now-unused fourth breakpoint value in column E.
$ws.Range("C9").Value = 1.7
$ws.Range("D9").Value = 2
$ws.Range("E9").ClearContents()

# Row 10 ("nav1"): shrink the base font size.
$ws.Range("B10").Value = 0.95

# Leave the selection where the author ended up editing.
$ws.Range("B10").Select()
